$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 16.465
$ws.Range("C3").Value = -13.008
$ws.Range("C14").Value = -12.139
$ws.Range("C16").Value = -12.893
$ws.Range("E18").Value = 17.148
$ws.Range("C21").Value = -12.842
$ws.Range("C23").Value = -12.223
$ws.Range("E24").Value = 17.003
$ws.Range("C25").Value = -11.549
$ws.Range("E25").Value = 16.98
$ws.Range("C26").Value = -12.33
$ws.Range("E27").Value = 16.147
$ws.Range("C29").Value = -12.061
$ws.Range("E30").Value = 16.358
$ws.Range("E31").Value = 16.411
$ws.Range("E39").Value = 16.508
$ws.Range("C40").Value = -12.426
$ws.Range("E42").Value = 16.663
$ws.Range("E48").Value = 17.232
$ws.Range("E51").Value = 16.878
$ws.Range("E52").Value = 16.954
$ws.Range("C53").Value = -11.515
$ws.Range("E55").Value = 16.491
$ws.Range("E56").Value = 16.414
$ws.Range("C57").Value = -13.794
$ws.Range("E57").Value = 16.572
$ws.Range("C59").Value = -13.155
$ws.Range("E60").Value = 16.592
$ws.Range("C65").Value = -12.118
$ws.Range("C69").Value = -10.697
$ws.Range("E73").Value = 16.778
$ws.Range("E74").Value = 16.319
$ws.Range("C79").Value = -12.083
$ws.Range("C83").Value = -13.069
$ws.Range("E89").Value = 17.158
$ws.Range("E90").Value = 16.483
$ws.Range("C91").Value = -11.133
$ws.Range("E92").Value = 17.779
$ws.Range("C93").Value = -11.624
$ws.Range("C100").Value = -12.727
